$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.528.57'
$ws.Range("E2").Value = '  -2.24%  '

$ws.Range("D3").Value = '1.961.83'
$ws.Range("E3").Value = '  -3.91%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '251.16'
$ws.Range("E5").Value = '  +1.34%  '

$ws.Range("D6").Value = '0.605'
$ws.Range("E6").Value = '  -3.18%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '54.17'
$ws.Range("E8").Value = '  -8.74%  '

$ws.Range("D9").Value = '0.370'
$ws.Range("E9").Value = '  -5.97%  '

$ws.Range("D10").Value = '0.0749'
$ws.Range("E10").Value = '  -7.24%  '

$ws.Range("E11").Value = '  -3.75%  '

$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '13.97'
$ws.Range("E12").Value = '  -8.30%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.253.19'
$ws.Range("E13").Value = '  -3.83%  '

$ws.Range("D14").Value = '20.99'
$ws.Range("E14").Value = '  -5.05%  '

$ws.Range("D15").Value = '0.766'
$ws.Range("E15").Value = '  -10.10%  '

$ws.Range("D16").Value = '5.10'
$ws.Range("E16").Value = '  -6.57%  '

$ws.Range("D17").Value = '1.959.08'
$ws.Range("E17").Value = '  -3.95%  '

$ws.Range("D18").Value = '36.371.11'
$ws.Range("E18").Value = '  -2.60%  '

$ws.Range("D19").Value = '68.88'
$ws.Range("E19").Value = '  -2.29%  '

$ws.Range("D20").Value = '0.0₃0814'
$ws.Range("E20").Value = '  -5.69%  '

$ws.Range("D21").Value = '230.32'
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").Value = '4.98'
$ws.Range("E22").Value = '  -5.46%  '

$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").Value = '2.49'
$ws.Range("E24").Value = '  -2.60%  '

$ws.Range("E25").Value = '  -0.12%  '

$ws.Range("D26").Value = '162.76'
$ws.Range("E26").Value = '  -1.11%  '

$ws.Range("D27").Value = '8.73'
$ws.Range("E27").Value = '  -7.18%  '

$ws.Range("D28").Value = '19.07'
$ws.Range("E28").Value = '  -4.32%  '

$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = '1.32'
$ws.Range("E29").Value = '  -4.02%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '0.121'
$ws.Range("E30").Value = '  -12.06%  '

$ws.Range("E31").Value = '  -3.68%  '

$ws.Range("D32").Value = '4.45'
$ws.Range("E32").Value = '  -7.00%  '

$ws.Range("E33").Value = '  -9.25%  '

$ws.Range("D34").Value = '4.27'
$ws.Range("E34").Value = '  -5.46%  '

$ws.Range("D35").Value = '2.32'
$ws.Range("E35").Value = '  -9.43%  '

$ws.Range("D36").Value = '1.80'
$ws.Range("E36").Value = '  -0.27%  '

$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.27%  '

$ws.Range("D38").Value = '3.32'
$ws.Range("E38").Value = '  -6.67%  '

$ws.Range("D39").Value = '5.30'
$ws.Range("E39").Value = '  -3.05%  '

$ws.Range("D40").Value = '2.97'
$ws.Range("E40").Value = '  -0.94%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.436.65'
$ws.Range("E41").Value = '  +3.68%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.15'
$ws.Range("E42").Value = '  -3.15%  '

$ws.Range("D43").Value = '0.0899'
$ws.Range("E43").Value = '  -8.35%  '

$ws.Range("D44").Value = '0.0204'
$ws.Range("E44").Value = '  -5.77%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '15.30'
$ws.Range("E45").Value = '  -8.43%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '86.92'
$ws.Range("E46").Value = '  -5.40%  '

$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -6.51%  '

$ws.Range("D48").Value = '2.86'
$ws.Range("E48").Value = '  -1.12%  '

$ws.Range("E49").Value = '  -10.56%  '

$ws.Range("D50").Value = '2.149.79'
$ws.Range("E50").Value = '  -3.81%  '

$ws.Range("D51").Value = '1.90'
$ws.Range("E51").Value = '  -10.50%  '
